$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 312865.88
$ws.Range("I38").Value = 434889.7
$ws.Range("J38").Value = 1027.2222
$ws.Range("K38").Value = 1304669.1
$ws.Range("L38").Value = 3081.6666
$ws.Range("M38").Value = -1304297.1
$ws.Range("N38").Value = -3825.6666
$ws.Range("H39").Value = 1892.6923
$ws.Range("I39").Value = 84
$ws.Range("J39").Value = 4002.8333
$ws.Range("K39").Value = 252
$ws.Range("L39").Value = 12008.4999
$ws.Range("M39").Value = 44
$ws.Range("N39").Value = -12600.4999
$ws.Range("H88").Value = 1711.5555
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 2601
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 2601
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -3413
$ws.Range("H91").Value = 1711.5555
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 2601
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 2601
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -5409
$ws.Range("H124").Value = 60780
$ws.Range("J124").Value = 60780
$ws.Range("L124").Value = 60780
$ws.Range("N124").Value = -70600
$ws.Range("H135").Value = 26316592
$ws.Range("I135").Value = 421.5
$ws.Range("J135").Value = 100001864
$ws.Range("K135").Value = 3793.5
$ws.Range("L135").Value = 900016776
$ws.Range("M135").Value = -1258.5
$ws.Range("N135").Value = -900021846
$ws.Range("H141").Value = 2802.8845
$ws.Range("I141").Value = 2515
$ws.Range("J141").Value = 3450.625
$ws.Range("K141").Value = 7545
$ws.Range("L141").Value = 10351.875
$ws.Range("M141").Value = -2365
$ws.Range("N141").Value = -20711.875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 129.27272
$ws.Range("I5").Value = 90
$ws.Range("J5").Value = 234
$ws.Range("K5").Value = 90
$ws.Range("L5").Value = 234
$ws.Range("M5").Value = 22
$ws.Range("N5").Value = -458
$ws.Range("H35").Value = 42020.5
$ws.Range("I35").Value = 4000
$ws.Range("K35").Value = 4000
$ws.Range("M35").Value = -3594
$ws.Range("H128").Value = 55000
$ws.Range("J128").Value = 55000
$ws.Range("L128").Value = 55000
$ws.Range("N128").Value = -64960

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 129.27272
$ws.Range("I4").Value = 90
$ws.Range("J4").Value = 234
$ws.Range("K4").Value = 90
$ws.Range("L4").Value = 234
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = -464

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1299.8
$ws.Range("I22").Value = 342.2857
$ws.Range("J22").Value = 3534
$ws.Range("K22").Value = 342.2857
$ws.Range("L22").Value = 3534
$ws.Range("M22").Value = 7.71429999999998
$ws.Range("N22").Value = -4234
$ws.Range("H59").Value = 8488.888999999999
$ws.Range("J59").Value = 8488.888999999999
$ws.Range("L59").Value = 8488.888999999999
$ws.Range("N59").Value = -10778.889
$ws.Range("H132").Value = 2402.4707
$ws.Range("I132").Value = 1653.6666
$ws.Range("J132").Value = 4199.6
$ws.Range("K132").Value = 4960.9998
$ws.Range("L132").Value = 12598.8
$ws.Range("M132").Value = -2430.9998
$ws.Range("N132").Value = -17658.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 198.53334
$ws.Range("I14").Value = 198.53334
$ws.Range("K14").Value = 595.6000200000001
$ws.Range("M14").Value = -422.6000200000001
$ws.Range("H34").Value = 412.80768
$ws.Range("I34").Value = 92
$ws.Range("K34").Value = 276
$ws.Range("M34").Value = -192
$ws.Range("H39").Value = 13864
$ws.Range("I39").Value = 700
$ws.Range("J39").Value = 15180.4
$ws.Range("K39").Value = 2100
$ws.Range("L39").Value = 45541.2
$ws.Range("M39").Value = -1806
$ws.Range("N39").Value = -46129.2
$ws.Range("H55").Value = 78577096
$ws.Range("I55").Value = 4250
$ws.Range("J55").Value = 91672570
$ws.Range("K55").Value = 12750
$ws.Range("L55").Value = 275017710
$ws.Range("N55").Value = -275018064
$ws.Range("M55").Value = -12573
$ws.Range("H131").Value = 1071.8
$ws.Range("J131").Value = 1152.0769
$ws.Range("L131").Value = 3456.2307
$ws.Range("N131").Value = -13536.2307

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1750.75
$ws.Range("J22").Value = 1901
$ws.Range("L22").Value = 1901
$ws.Range("N22").Value = -2491
$ws.Range("H27").Value = 1750.75
$ws.Range("J27").Value = 1901
$ws.Range("L27").Value = 1901
$ws.Range("N27").Value = -2115
$ws.Range("H68").Value = 2064.158
$ws.Range("I68").Value = 1818.091
$ws.Range("J68").Value = 2402.5
$ws.Range("K68").Value = 1818.091
$ws.Range("L68").Value = 2402.5
$ws.Range("M68").Value = -1069.091
$ws.Range("N68").Value = -3900.5
$ws.Range("H71").Value = 2064.158
$ws.Range("I71").Value = 1818.091
$ws.Range("J71").Value = 2402.5
$ws.Range("K71").Value = 9090.455
$ws.Range("L71").Value = 12012.5
$ws.Range("M71").Value = -5346.455
$ws.Range("N71").Value = -19500.5
$ws.Range("H132").Value = 461701.97
$ws.Range("I132").Value = 129930.5
$ws.Range("J132").Value = 627587.7
$ws.Range("K132").Value = 389791.5
$ws.Range("L132").Value = 1882763.1
$ws.Range("M132").Value = -387261.5
$ws.Range("N132").Value = -1887823.1
$ws.Range("H136").Value = 264175.2
$ws.Range("I136").Value = 400483.1
$ws.Range("J136").Value = 2044.6154
$ws.Range("K136").Value = 1201449.3
$ws.Range("L136").Value = 6133.8462
$ws.Range("M136").Value = -1198899.3
$ws.Range("N136").Value = -11233.8462

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 33950
$ws.Range("J123").Value = 33950
$ws.Range("L123").Value = 33950
$ws.Range("N123").Value = -43750
$ws.Range("H132").Value = 3838.8223
$ws.Range("I132").Value = 1357.3871
$ws.Range("J132").Value = 9333.429
$ws.Range("K132").Value = 4072.1613
$ws.Range("L132").Value = 28000.287
$ws.Range("M132").Value = -1542.1613
$ws.Range("N132").Value = -33060.287
$ws.Range("H136").Value = 264582.16
$ws.Range("I136").Value = 1656.6086
$ws.Range("K136").Value = 4969.825800000001
$ws.Range("M136").Value = -2419.825800000001
